$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.380.52'
Set-TextValue $ws.Range('E2') '  -0.09%  '
Set-TextValue $ws.Range('D3') '1.846.74'
Set-TextValue $ws.Range('E3') '  -0.02%  '
Set-TextValue $ws.Range('D4') '0.9994'
Set-TextValue $ws.Range('E4') '  +0.10%  '
Set-TextValue $ws.Range('D5') '240.03'
Set-TextValue $ws.Range('D6') '0.6297'
Set-TextValue $ws.Range('E6') '  -0.56%  '
Set-TextValue $ws.Range('E7') '  +0.02%  '
Set-TextValue $ws.Range('D8') '0.07600'
Set-TextValue $ws.Range('E8') '  +0.56%  '
Set-TextValue $ws.Range('D9') '0.2930'
Set-TextValue $ws.Range('E9') '  -1.24%  '
Set-TextValue $ws.Range('D10') '24.50'
Set-TextValue $ws.Range('E10') '  -0.56%  '
Set-TextValue $ws.Range('D11') '0.07742'
Set-TextValue $ws.Range('E11') '  +0.13%  '
Set-TextValue $ws.Range('D12') '1.840.65'
Set-TextValue $ws.Range('E12') '  -0.47%  '
Set-TextValue $ws.Range('D13') '0.00001091'
Set-TextValue $ws.Range('E13') '  +9.38%  '
Set-TextValue $ws.Range('D14') '5.001'
Set-TextValue $ws.Range('E14') '  +0.00%  '
Set-TextValue $ws.Range('D15') '0.6779'
Set-TextValue $ws.Range('E15') '  -1.05%  '
Set-TextValue $ws.Range('D16') '83.65'
Set-TextValue $ws.Range('E16') '  +0.70%  '
Set-TextValue $ws.Range('D17') '2.089.34'
Set-TextValue $ws.Range('E17') '  -7.73%  '
Set-TextValue $ws.Range('D18') '6.155'
Set-TextValue $ws.Range('E18') '  -0.39%  '
Set-TextValue $ws.Range('D19') '29.407.99'
Set-TextValue $ws.Range('E19') '  -0.06%  '
Set-TextValue $ws.Range('D20') '228.68'
Set-TextValue $ws.Range('E20') '  -0.60%  '
Set-TextValue $ws.Range('D21') '12.43'
Set-TextValue $ws.Range('E21') '  -0.42%  '
Set-TextValue $ws.Range('E22') '  +0.05%  '
Set-TextValue $ws.Range('D23') '7.419'
Set-TextValue $ws.Range('E23') '  -2.04%  '
Set-TextValue $ws.Range('D25') '157.00'
Set-TextValue $ws.Range('E25') '  -0.09%  '
Set-TextValue $ws.Range('D26') '0.1393'
Set-TextValue $ws.Range('E26') '  -0.79%  '
Set-TextValue $ws.Range('D27') '8.384'
Set-TextValue $ws.Range('E27') '  +0.00%  '
Set-TextValue $ws.Range('E28') '  -0.36%  '
Set-TextValue $ws.Range('D29') '1.465'
Set-TextValue $ws.Range('E29') '  +0.04%  '
Set-TextValue $ws.Range('D30') '1.309'
Set-TextValue $ws.Range('E30') '  +4.58%  '
Set-TextValue $ws.Range('D31') '0.05616'
Set-TextValue $ws.Range('E31') '  -1.97%  '
Set-TextValue $ws.Range('D32') '4.101'
Set-TextValue $ws.Range('E32') '  -0.71%  '
Set-TextValue $ws.Range('D33') '4.036'
Set-TextValue $ws.Range('E33') '  -0.02%  '
Set-TextValue $ws.Range('D34') '1.846'
Set-TextValue $ws.Range('E34') '  -0.43%  '
Set-TextValue $ws.Range('E35') '  -0.03%  '
Set-TextValue $ws.Range('D36') '0.7102'
Set-TextValue $ws.Range('E36') '  -0.90%  '
Set-TextValue $ws.Range('D37') '2.582'
Set-TextValue $ws.Range('E37') '  -0.49%  '
Set-TextValue $ws.Range('D38') '1.231.85'
Set-TextValue $ws.Range('E38') '  -1.62%  '
Set-TextValue $ws.Range('E39') '  -0.62%  '
Set-TextValue $ws.Range('D40') '2.766'
Set-TextValue $ws.Range('E40') '  -0.56%  '
Set-TextValue $ws.Range('D41') '6.446'
Set-TextValue $ws.Range('E41') '  +3.80%  '
Set-TextValue $ws.Range('D42') '0.9071'
Set-TextValue $ws.Range('E42') '  -0.16%  '
Set-TextValue $ws.Range('D43') '0.9998'
Set-TextValue $ws.Range('E43') '  -0.05%  '
Set-TextValue $ws.Range('B44') 'RocketPoolETH'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range('D44') '1.998.88'
Set-TextValue $ws.Range('E44') '  -0.64%  '
Set-TextValue $ws.Range('B45') 'Quant'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D45') '101.58'
Set-TextValue $ws.Range('E45') '  -0.21%  '
Set-TextValue $ws.Range('B46') 'Aave'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D46') '66.04'
Set-TextValue $ws.Range('E46') '  -0.77%  '
Set-TextValue $ws.Range('B47') 'BabyDogeCoin'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D47') '0.00000000123'
Set-TextValue $ws.Range('E47') '  +5.59%  '
Set-TextValue $ws.Range('B48') 'Aptos'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D48') '7.188'
Set-TextValue $ws.Range('E48') '  +1.35%  '
Set-TextValue $ws.Range('B49') 'TheSandbox'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D49') '0.4016'
Set-TextValue $ws.Range('E49') '  -0.33%  '
Set-TextValue $ws.Range('B50') 'EnergySwap'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D50') '8.995'
Set-TextValue $ws.Range('E50') '  -2.08%  '
Set-TextValue $ws.Range('B51') 'RenderToken'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D51') '1.680'
Set-TextValue $ws.Range('E51') '  -1.99%  '
